$d = $word.ActiveDocument

$p = $d.Paragraphs.Last

$p.Range.ParagraphFormat.LineSpacingRule = 0
$p.Range.ParagraphFormat.LineSpacing = 24

$p.Range.Font.NameAscii = "Arial"
$p.Range.Font.NameOther = "Arial"
$p.Range.Font.NameBi = "Arial"

Write-Output "Applied formatting to the trailing empty paragraph."
